$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.507.97'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '2.438.42'
$ws.Range("E3").Value = '  -0.28%  '
$ws.Range("E4").Value = '  -0.19%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.98'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.05%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  -0.73%  '
$ws.Range("D9").Value = '2.434.97'
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.81%  '
$ws.Range("E11").Value = '  +0.80%  '
$ws.Range("E12").Value = '  -0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.352'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.36%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.56'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("E15").Value = '  -2.94%  '
$ws.Range("D16").Value = '2.868.42'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = '62.480.38'
$ws.Range("E17").Value = '  -1.16%  '
$ws.Range("D18").Value = '2.444.87'
$ws.Range("E18").Value = '  +0.10%  '
$ws.Range("E19").Value = '  -2.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.14'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '325.94'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.93%  '
$ws.Range("E23").Value = '  +5.07%  '
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '65.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.80%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '626.86'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.60%  '
$ws.Range("D28").Value = '0.0₃0961'
$ws.Range("E28").Value = '  -6.76%  '
$ws.Range("D29").Value = '2.561.11'
$ws.Range("E29").Value = '  -1.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.970'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.99%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.43'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.10%  '
$ws.Range("E32").Value = '  -2.09%  '
$ws.Range("E33").Value = '  -1.33%  '
$ws.Range("E34").Value = '  -6.78%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.98'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.19%  '
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("E37").Value = '  -3.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.374'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.79%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.46'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '146.02'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '42.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.20%  '
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  -5.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '145.47'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0523'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.92%  '
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.65'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.70%  '
$ws.Range("E51").Value = '  -2.20%  '
